$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "user-name" column (B) so it becomes column C,
# and the new column B holds "client-id" with value 1130 for every data row.
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "client-id"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Value = 1130
}

$ws.Range("C11").Select()
